$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.566.71"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "2.486.66"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.48%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Value = "2.496.74"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0980"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "2.913.83"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "56.631.62"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.502.25"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").Value = "2.598.00"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  +3.81%  "
$ws.Range("D30").Value = "0.0₃0798"
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  +4.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0557"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.611"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "263.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.60%  "
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "1.901.89"
$ws.Range("E51").Value = "  -3.41%  "
